$p = $ppt.ActivePresentation

# --- Slide 2: "Money savings for you with hourly rate!" content placeholder ---
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

$tr2.Paragraphs(1,1).Runs(1,1).Text = "Cut costs by planning!"
$tr2.Paragraphs(2,1).Runs(1,1).Text = "Is it cheaper to charge your electric car tomorrow or the day after?"
$tr2.Paragraphs(3,1).Runs(1,1).Text = "Heat the house just before the weather turns cold and prices rise!"

# --- Slide 3: "Easy-to-use" content placeholder ---
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(4)

# Resize / reposition the placeholder
$shp3.Left = 61
$shp3.Top = 157.1501
$shp3.Width = 323.75
$shp3.Height = 371.6249

$tr3 = $shp3.TextFrame.TextRange

# Update last paragraph's text/size, then append a new paragraph
$tr3.Paragraphs(5,1).Runs(1,1).Text = "Just three recommendation levels per hour "
$tr3.Paragraphs(5,1).Runs(1,1).Font.Size = 26
$null = $tr3.InsertAfter([char]13 + "One pageload!")
$tr3.Paragraphs(6,1).Runs(1,1).Font.Size = 26
